{"js": "// Translate the three Spanish cover-page labels to English:\n//   CURSO:       -> COURSE:\n//   PROFESOR:    -> PROFESSOR:\n//   INTEGRANTES: -> MEMBERS:\n// Only the label word itself changes; the rest of each line (the colon,\n// the value after it, etc.) is left untouched, exactly like selecting the\n// word in Word and typing its replacement.\n\nconst body = context.document.body;\n\nasync function replaceWord(original, replacement) {\n  const results = body.search(original, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceWord(\"CURSO\", \"COURSE\");\nawait replaceWord(\"PROFESOR\", \"PROFESSOR\");\nawait replaceWord(\"INTEGRANTES\", \"MEMBERS\");\n", "ps1": "# Translate the three Spanish cover-page labels to English:\n#   CURSO:       -> COURSE:\n#   PROFESOR:    -> PROFESSOR:\n#   INTEGRANTES: -> MEMBERS:\n# Only the label word itself is replaced; everything after it on the same\n# line (the colon, the value, etc.) is left untouched.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Word($find, $replacement) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute(\n        $find,        # FindText\n        $true,        # MatchCase\n        $true,        # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replacement, # ReplaceWith\n        $wdReplaceAll # Replace\n    ) | Out-Null\n}\n\nReplace-Word \"CURSO\" \"COURSE\"\nReplace-Word \"PROFESOR\" \"PROFESSOR\"\nReplace-Word \"INTEGRANTES\" \"MEMBERS\"\n"}
